$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in previously-empty name_CN cells for existing rows 80 and 81
$ws.Range("E80").Value = "Ｅ－火瑪麗"
$ws.Range("E81").Value = "Ｅ－水瑪麗"

# Add two new rows (82, 83) for the JP tenth anniversary servant collections
$ws.Range("A82").Value = 436
$ws.Range("B82").Value = 0
$ws.Range("C82").Value = "Uolgamariegrandcollection"
$ws.Range("D82").Value = "Ｅ－グランマリー"
$ws.Range("E82").Value = ""
$ws.Range("F82").Value = ""

$ws.Range("A83").Value = 443
$ws.Range("B83").Value = 0
$ws.Range("C83").Value = "Uolgamariestellarcollection"
$ws.Range("D83").Value = "Ｅ－ステラマリー"
$ws.Range("E83").Value = ""
$ws.Range("F83").Value = ""
